$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain numeric-looking Price values to remain text (matches source formatting)
$textRows = @(4,5,6,7,8,9,10,11,12,14,15,16,17,18,19,20,21,23,24,25,26,27,28,29,30,31,32,33,34,35,37,38,39,40,41,42,44,45,46,47,48,49,50,51)
foreach ($r in $textRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# Update Price (D) and Volume(1h) (E) columns with latest values
$ws.Range("D2").Value = "27.609.08"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").Value = "1.851.34"
$ws.Range("E3").Value = "  -2.18%  "
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  -0.90%  "
$ws.Range("D5").Value = "334.48"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").Value = "1.009"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("D7").Value = "0.4601"
$ws.Range("E7").Value = "  -1.93%  "
$ws.Range("D8").Value = "0.3890"
$ws.Range("E8").Value = "  -0.82%  "
$ws.Range("D9").Value = "45.65"
$ws.Range("E9").Value = "  -4.03%  "
$ws.Range("D10").Value = "0.07929"
$ws.Range("E10").Value = "  -1.23%  "
$ws.Range("D11").Value = "1.002"
$ws.Range("E11").Value = "  -1.80%  "
$ws.Range("D12").Value = "21.55"
$ws.Range("E12").Value = "  -1.17%  "
$ws.Range("D13").Value = "1.856.49"
$ws.Range("E13").Value = "  -2.06%  "
$ws.Range("D14").Value = "5.949"
$ws.Range("D15").Value = "7.169"
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("D16").Value = "1.011"
$ws.Range("E16").Value = "  -0.86%  "
$ws.Range("D17").Value = "88.35"
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("D18").Value = "0.06704"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("D19").Value = "0.00001033"
$ws.Range("E19").Value = "  -1.73%  "
$ws.Range("D20").Value = "17.25"
$ws.Range("E20").Value = "  +0.50%  "
$ws.Range("D21").Value = "1.009"
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("D22").Value = "27.617.66"
$ws.Range("E22").Value = "  -1.27%  "
$ws.Range("D23").Value = "5.408"
$ws.Range("E23").Value = "  -1.92%  "
$ws.Range("D24").Value = "10.91"
$ws.Range("E24").Value = "  -0.66%  "
$ws.Range("D25").Value = "2.308"
$ws.Range("E25").Value = "  -1.64%  "
$ws.Range("D26").Value = "158.78"
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("D27").Value = "19.53"
$ws.Range("E27").Value = "  -2.68%  "
$ws.Range("D28").Value = "2.128"
$ws.Range("E28").Value = "  +2.36%  "
$ws.Range("D29").Value = "5.440"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").Value = "121.20"
$ws.Range("E30").Value = "  -0.53%  "
$ws.Range("D31").Value = "0.9734"
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("D32").Value = "0.09408"
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("D33").Value = "3.613"
$ws.Range("E33").Value = "  -1.61%  "
$ws.Range("D34").Value = "5.303"
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("D35").Value = "1.344"
$ws.Range("E35").Value = "  -3.87%  "
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("D37").Value = "0.06009"
$ws.Range("E37").Value = "  -1.92%  "
$ws.Range("D38").Value = "8.355"
$ws.Range("E38").Value = "  +3.12%  "
$ws.Range("D39").Value = "1.185"
$ws.Range("E39").Value = "  -2.86%  "
$ws.Range("D40").Value = "1.009"
$ws.Range("E40").Value = "  -0.59%  "
$ws.Range("D41").Value = "0.5923"
$ws.Range("E41").Value = "  -1.09%  "
$ws.Range("D42").Value = "10.39"
$ws.Range("E42").Value = "  +0.44%  "
$ws.Range("E43").Value = "  -1.34%  "
$ws.Range("D44").Value = "1.242"
$ws.Range("E44").Value = "  -2.18%  "
$ws.Range("D45").Value = "0.5583"
$ws.Range("E45").Value = "  -1.96%  "
$ws.Range("D46").Value = "12.12"
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("D47").Value = "1.908"
$ws.Range("E47").Value = "  -1.49%  "
$ws.Range("D48").Value = "0.06706"
$ws.Range("E48").Value = "  -3.33%  "
$ws.Range("D49").Value = "111.18"
$ws.Range("E49").Value = "  -2.39%  "
$ws.Range("D50").Value = "1.050"
$ws.Range("E50").Value = "  -1.98%  "
$ws.Range("D51").Value = "1.010"
$ws.Range("E51").Value = "  -0.79%  "
